$d = $word.ActiveDocument
$d.Content.Find.Execute("43-34=", $true, $false, $false, $false, $false, $true, 1, $false, "66-58=", 2) | Out-Null
$d.Content.Find.Execute("64-21=", $true, $false, $false, $false, $false, $true, 1, $false, "87-45=", 2) | Out-Null
$d.Content.Find.Execute("24+30=", $true, $false, $false, $false, $false, $true, 1, $false, "19+25=", 2) | Out-Null
$d.Content.Find.Execute("42+16=", $true, $false, $false, $false, $false, $true, 1, $false, "90-84=", 2) | Out-Null
$d.Content.Find.Execute("10+54=", $true, $false, $false, $false, $false, $true, 1, $false, "32+37=", 2) | Out-Null
$d.Content.Find.Execute("26+27=", $true, $false, $false, $false, $false, $true, 1, $false, "55-33=", 2) | Out-Null
$d.Content.Find.Execute("95-50=", $true, $false, $false, $false, $false, $true, 1, $false, "17+80=", 2) | Out-Null
$d.Content.Find.Execute("39-1=", $true, $false, $false, $false, $false, $true, 1, $false, "43-30=", 2) | Out-Null
$d.Content.Find.Execute("46+12=", $true, $false, $false, $false, $false, $true, 1, $false, "5+71=", 2) | Out-Null
$d.Content.Find.Execute("13+75=", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=", 2) | Out-Null
$d.Content.Find.Execute("4+29=", $true, $false, $false, $false, $false, $true, 1, $false, "72+4=", 2) | Out-Null
$d.Content.Find.Execute("79-50=", $true, $false, $false, $false, $false, $true, 1, $false, "32+17=", 2) | Out-Null
$d.Content.Find.Execute("14+19=", $true, $false, $false, $false, $false, $true, 1, $false, "19+61=", 2) | Out-Null
$d.Content.Find.Execute("61-26=", $true, $false, $false, $false, $false, $true, 1, $false, "12+67=", 2) | Out-Null
$d.Content.Find.Execute("26-14=", $true, $false, $false, $false, $false, $true, 1, $false, "79-35=", 2) | Out-Null
$d.Content.Find.Execute("56+43=", $true, $false, $false, $false, $false, $true, 1, $false, "54+15=", 2) | Out-Null
$d.Content.Find.Execute("21+48=", $true, $false, $false, $false, $false, $true, 1, $false, "96-69=", 2) | Out-Null
$d.Content.Find.Execute("60+21=", $true, $false, $false, $false, $false, $true, 1, $false, "27+67=", 2) | Out-Null
$d.Content.Find.Execute("72-61=", $true, $false, $false, $false, $false, $true, 1, $false, "94-80=", 2) | Out-Null
$d.Content.Find.Execute("27+57=", $true, $false, $false, $false, $false, $true, 1, $false, "23-19=", 2) | Out-Null
$d.Content.Find.Execute("28+10=", $true, $false, $false, $false, $false, $true, 1, $false, "45-33=", 2) | Out-Null
$d.Content.Find.Execute("68-18=", $true, $false, $false, $false, $false, $true, 1, $false, "47-42=", 2) | Out-Null
$d.Content.Find.Execute("95-91=", $true, $false, $false, $false, $false, $true, 1, $false, "33+65=", 2) | Out-Null
$d.Content.Find.Execute("65-34=", $true, $false, $false, $false, $false, $true, 1, $false, "61-23=", 2) | Out-Null
$d.Content.Find.Execute("22+41=", $true, $false, $false, $false, $false, $true, 1, $false, "71-44=", 2) | Out-Null
$d.Content.Find.Execute("2+80=", $true, $false, $false, $false, $false, $true, 1, $false, "1+66=", 2) | Out-Null
$d.Content.Find.Execute("48+2=", $true, $false, $false, $false, $false, $true, 1, $false, "56+9=", 2) | Out-Null
$d.Content.Find.Execute("60+31=", $true, $false, $false, $false, $false, $true, 1, $false, "96-68=", 2) | Out-Null
$d.Content.Find.Execute("12+65=", $true, $false, $false, $false, $false, $true, 1, $false, "89-42=", 2) | Out-Null
$d.Content.Find.Execute("85-19=", $true, $false, $false, $false, $false, $true, 1, $false, "11+43=", 2) | Out-Null
$d.Content.Find.Execute("87-42=", $true, $false, $false, $false, $false, $true, 1, $false, "41+7=", 2) | Out-Null
$d.Content.Find.Execute("39+15=", $true, $false, $false, $false, $false, $true, 1, $false, "75-34=", 2) | Out-Null
$d.Content.Find.Execute("38+40=", $true, $false, $false, $false, $false, $true, 1, $false, "58+31=", 2) | Out-Null
$d.Content.Find.Execute("42-25=", $true, $false, $false, $false, $false, $true, 1, $false, "38+20=", 2) | Out-Null
$d.Content.Find.Execute("83-28=", $true, $false, $false, $false, $false, $true, 1, $false, "96-94=", 2) | Out-Null
$d.Content.Find.Execute("12-2=", $true, $false, $false, $false, $false, $true, 1, $false, "92-53=", 2) | Out-Null
$d.Content.Find.Execute("60-23=", $true, $false, $false, $false, $false, $true, 1, $false, "28-6=", 2) | Out-Null
$d.Content.Find.Execute("23+2=", $true, $false, $false, $false, $false, $true, 1, $false, "50+17=", 2) | Out-Null
$d.Content.Find.Execute("64+20=", $true, $false, $false, $false, $false, $true, 1, $false, "72-58=", 2) | Out-Null
$d.Content.Find.Execute("90-26=", $true, $false, $false, $false, $false, $true, 1, $false, "25+40=", 2) | Out-Null
$d.Content.Find.Execute("46+25=", $true, $false, $false, $false, $false, $true, 1, $false, "5+45=", 2) | Out-Null
$d.Content.Find.Execute("68-0=", $true, $false, $false, $false, $false, $true, 1, $false, "57-40=", 2) | Out-Null
$d.Content.Find.Execute("79-71=", $true, $false, $false, $false, $false, $true, 1, $false, "46+6=", 2) | Out-Null
$d.Content.Find.Execute("17+11=", $true, $false, $false, $false, $false, $true, 1, $false, "83+9=", 2) | Out-Null
$d.Content.Find.Execute("12-0=", $true, $false, $false, $false, $false, $true, 1, $false, "8+35=", 2) | Out-Null
$d.Content.Find.Execute("24-22=", $true, $false, $false, $false, $false, $true, 1, $false, "4+94=", 2) | Out-Null
$d.Content.Find.Execute("86-66=", $true, $false, $false, $false, $false, $true, 1, $false, "82-15=", 2) | Out-Null
$d.Content.Find.Execute("78+9=", $true, $false, $false, $false, $false, $true, 1, $false, "40+43=", 2) | Out-Null
$d.Content.Find.Execute("36+36=", $true, $false, $false, $false, $false, $true, 1, $false, "18+70=", 2) | Out-Null
$d.Content.Find.Execute("8+14=", $true, $false, $false, $false, $false, $true, 1, $false, "52+25=", 2) | Out-Null
$d.Content.Find.Execute("1+36=", $true, $false, $false, $false, $false, $true, 1, $false, "95-54=", 2) | Out-Null
$d.Content.Find.Execute("59-39=", $true, $false, $false, $false, $false, $true, 1, $false, "91-0=", 2) | Out-Null
$d.Content.Find.Execute("47+13=", $true, $false, $false, $false, $false, $true, 1, $false, "18+39=", 2) | Out-Null
$d.Content.Find.Execute("67-27=", $true, $false, $false, $false, $false, $true, 1, $false, "48-21=", 2) | Out-Null
$d.Content.Find.Execute("26-13=", $true, $false, $false, $false, $false, $true, 1, $false, "72-26=", 2) | Out-Null
$d.Content.Find.Execute("71-16=", $true, $false, $false, $false, $false, $true, 1, $false, "64+31=", 2) | Out-Null
$d.Content.Find.Execute("15+46=", $true, $false, $false, $false, $false, $true, 1, $false, "81-11=", 2) | Out-Null
$d.Content.Find.Execute("78-8=", $true, $false, $false, $false, $false, $true, 1, $false, "89-28=", 2) | Out-Null
$d.Content.Find.Execute("27+61=", $true, $false, $false, $false, $false, $true, 1, $false, "95+4=", 2) | Out-Null
$d.Content.Find.Execute("53-9=", $true, $false, $false, $false, $false, $true, 1, $false, "3+60=", 2) | Out-Null
$d.Content.Find.Execute("9+18=", $true, $false, $false, $false, $false, $true, 1, $false, "89-18=", 2) | Out-Null
$d.Content.Find.Execute("87+2=", $true, $false, $false, $false, $false, $true, 1, $false, "43-32=", 2) | Out-Null
$d.Content.Find.Execute("73+1=", $true, $false, $false, $false, $false, $true, 1, $false, "54-42=", 2) | Out-Null
$d.Content.Find.Execute("5+74=", $true, $false, $false, $false, $false, $true, 1, $false, "33+58=", 2) | Out-Null
$d.Content.Find.Execute("40-12=", $true, $false, $false, $false, $false, $true, 1, $false, "65+33=", 2) | Out-Null
$d.Content.Find.Execute("33-5=", $true, $false, $false, $false, $false, $true, 1, $false, "32-24=", 2) | Out-Null
$d.Content.Find.Execute("32+20=", $true, $false, $false, $false, $false, $true, 1, $false, "46-17=", 2) | Out-Null
$d.Content.Find.Execute("51+25=", $true, $false, $false, $false, $false, $true, 1, $false, "87-83=", 2) | Out-Null
$d.Content.Find.Execute("31+49=", $true, $false, $false, $false, $false, $true, 1, $false, "35+44=", 2) | Out-Null
$d.Content.Find.Execute("99-89=", $true, $false, $false, $false, $false, $true, 1, $false, "8+12=", 2) | Out-Null
$d.Content.Find.Execute("6+24=", $true, $false, $false, $false, $false, $true, 1, $false, "49+3=", 2) | Out-Null
$d.Content.Find.Execute("79-78=", $true, $false, $false, $false, $false, $true, 1, $false, "42+50=", 2) | Out-Null
$d.Content.Find.Execute("40-21=", $true, $false, $false, $false, $false, $true, 1, $false, "12+76=", 2) | Out-Null
$d.Content.Find.Execute("55-46=", $true, $false, $false, $false, $false, $true, 1, $false, "98-32=", 2) | Out-Null
$d.Content.Find.Execute("12+82=", $true, $false, $false, $false, $false, $true, 1, $false, "10+58=", 2) | Out-Null
$d.Content.Find.Execute("21+16=", $true, $false, $false, $false, $false, $true, 1, $false, "88-46=", 2) | Out-Null
$d.Content.Find.Execute("83-21=", $true, $false, $false, $false, $false, $true, 1, $false, "68-57=", 2) | Out-Null
$d.Content.Find.Execute("86-38=", $true, $false, $false, $false, $false, $true, 1, $false, "4+72=", 2) | Out-Null
$d.Content.Find.Execute("70-56=", $true, $false, $false, $false, $false, $true, 1, $false, "80-43=", 2) | Out-Null
$d.Content.Find.Execute("96-88=", $true, $false, $false, $false, $false, $true, 1, $false, "60-43=", 2) | Out-Null
$d.Content.Find.Execute("53-32=", $true, $false, $false, $false, $false, $true, 1, $false, "30+10=", 2) | Out-Null
$d.Content.Find.Execute("36+49=", $true, $false, $false, $false, $false, $true, 1, $false, "44-20=", 2) | Out-Null
$d.Content.Find.Execute("69-1=", $true, $false, $false, $false, $false, $true, 1, $false, "95-43=", 2) | Out-Null
$d.Content.Find.Execute("85-72=", $true, $false, $false, $false, $false, $true, 1, $false, "56-34=", 2) | Out-Null
$d.Content.Find.Execute("75-61=", $true, $false, $false, $false, $false, $true, 1, $false, "85-56=", 2) | Out-Null
$d.Content.Find.Execute("76-50=", $true, $false, $false, $false, $false, $true, 1, $false, "75-40=", 2) | Out-Null
$d.Content.Find.Execute("39-29=", $true, $false, $false, $false, $false, $true, 1, $false, "43+18=", 2) | Out-Null
$d.Content.Find.Execute("7+47=", $true, $false, $false, $false, $false, $true, 1, $false, "57-31=", 2) | Out-Null
$d.Content.Find.Execute("2+52=", $true, $false, $false, $false, $false, $true, 1, $false, "38+54=", 2) | Out-Null
$d.Content.Find.Execute("93-45=", $true, $false, $false, $false, $false, $true, 1, $false, "34-15=", 2) | Out-Null
$d.Content.Find.Execute("88-54=", $true, $false, $false, $false, $false, $true, 1, $false, "59-35=", 2) | Out-Null
$d.Content.Find.Execute("57+26=", $true, $false, $false, $false, $false, $true, 1, $false, "54+33=", 2) | Out-Null
$d.Content.Find.Execute("27+71=", $true, $false, $false, $false, $false, $true, 1, $false, "5+5=", 2) | Out-Null
$d.Content.Find.Execute("12+12=", $true, $false, $false, $false, $false, $true, 1, $false, "71-7=", 2) | Out-Null
$d.Content.Find.Execute("66+16=", $true, $false, $false, $false, $false, $true, 1, $false, "4+19=", 2) | Out-Null
$d.Content.Find.Execute("70-41=", $true, $false, $false, $false, $false, $true, 1, $false, "53-48=", 2) | Out-Null
$d.Content.Find.Execute("75-10=", $true, $false, $false, $false, $false, $true, 1, $false, "40-9=", 2) | Out-Null
$d.Content.Find.Execute("30+63=", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=", 2) | Out-Null
$d.Content.Find.Execute("42-4=", $true, $false, $false, $false, $false, $true, 1, $false, "15+75=", 2) | Out-Null
$d.Content.Find.Execute("94-21=", $true, $false, $false, $false, $false, $true, 1, $false, "93-90=", 2) | Out-Null
